$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# row => @{ col = newValue }
$changes = @{
    2  = @{ E = 58;  F = 40;  H = 52 }
    15 = @{ E = 174; F = 98;  H = 139 }
    17 = @{ E = 136 }
    19 = @{ E = 69 }
    20 = @{ E = 7;   F = 3;   H = 6 }
    21 = @{ E = 3;   F = 2;   H = 2 }
    23 = @{ E = 8;   F = 4;   H = 6 }
    27 = @{ F = 8;   H = 12 }
    29 = @{ E = 20 }
    35 = @{ E = 13 }
    36 = @{ E = 116 }
    37 = @{ E = 61 }
    38 = @{ E = 85 }
    40 = @{ E = 27;  F = 18;  H = 21 }
    41 = @{ E = 49 }
    45 = @{ E = 28 }
    47 = @{ E = 63;  G = 10;  H = 50 }
    48 = @{ E = 39;  F = 25;  G = 6;  H = 31 }
    50 = @{ E = 31;  F = 12;  H = 20 }
    51 = @{ E = 14 }
    53 = @{ E = 7 }
    57 = @{ E = 18 }
    61 = @{ E = 35;  F = 15;  H = 25 }
    62 = @{ E = 51 }
    68 = @{ E = 19 }
    69 = @{ E = 18 }
    70 = @{ E = 48 }
    75 = @{ E = 19 }
    76 = @{ E = 56;  F = 22;  H = 39 }
    77 = @{ E = 62 }
    79 = @{ E = 44 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
